$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Note" column (D) to the table, mirroring the formatting of
# column C (Address(es)) for every row of the table.
$ws.Range("C2:C11").Copy()
$ws.Range("D2:D11").PasteSpecial(-4122)

# Match the new column's width to column C's width.
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(3).ColumnWidth()

# Header for the new column.
$ws.Range("D2").Value = "Note"

# Only the TMP112 @ 0x48 row gets a note.
$ws.Range("D8").Value = "DRS4 Temperature Sensor"
